$wb = $excel.ActiveWorkbook

# "Student" sheet: students now just click a link to continue as a student
# instead of logging in with credentials.
$wsStudent = $wb.Worksheets.Item("Student")
$wsStudent.Range("B4").Value2 = "Se apasă link-ul ""Continuă ca student""."

# "Admin" sheet: the list of visible users no longer includes students,
# since students no longer have to log in.
$wsAdmin = $wb.Worksheets.Item("Admin")
$newAdminText = "Adminul este întâmpinat de pagina de unde se pot vizualiza utilizatorii curenți ai aplicației (profesori)."
$wsAdmin.Range("C44").Value2 = $newAdminText
$wsAdmin.Range("C53").Value2 = $newAdminText

# Leave the selection/cursor on each sheet where the author's editing
# session left it.
$wsDidactic = $wb.Worksheets.Item("Cadru Didactic")
$wsDidactic.Select()
$wsDidactic.Range("G12").Select()

$wsAdmin.Select()
$wsAdmin.Range("F49").Select()

$wsStudent.Select()
$wsStudent.Range("F12").Select()
